# This script rewrites the weekly price-listing rows (2-36) of the single
# worksheet so that each row's Fecha/Calidad/Volumen/Precio*/Unidad/Origen/Kg
# fields match the values from the re-shuffled weekly source extract.
# Data rows are represented as an array of hashtables: one entry per
# worksheet row, columns D,I,J,K,L,M,N,O,P,Q in that order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @(44595, "Primera", 200, 600, 700, 650, "$/paquete 6 unidades", "Región Metropolitana", 108, 6),
    @(44230, "Primera", 100, 600, 700, 650, "$/paquete 6 unidades", "Región de Ñuble", 108, 6),
    @(44230, "Segunda", 50, 500, 500, 500, "$/paquete 6 unidades", "Región de Ñuble", 83, 6),
    @(44321, "Primera", 100, 600, 700, 650, "$/paquete 6 unidades", "Región de Ñuble", 108, 6),
    @(44321, "Segunda", 50, 500, 500, 500, "$/paquete 6 unidades", "Región de Ñuble", 83, 6),
    @(44657, "Primera", 200, 8000, 9000, 8500, "$/paquete 36 unidades", "Región Metropolitana", 236, 36),
    @(44665, "Primera", 200, 600, 700, 650, "$/paquete 6 unidades", "Región de Ñuble", 108, 6),
    @(44665, "Segunda", 100, 500, 500, 500, "$/paquete 6 unidades", "Región de Ñuble", 83, 6),
    @(44308, "Primera", 200, 600, 700, 650, "$/paquete 6 unidades", "Región de Ñuble", 108, 6),
    @(44308, "Segunda", 100, 500, 500, 500, "$/paquete 6 unidades", "Región de Ñuble", 83, 6),
    @(44293, "Primera", 100, 600, 700, 650, "$/paquete 6 unidades", "Región de Ñuble", 108, 6),
    @(44293, "Segunda", 50, 500, 500, 500, "$/paquete 6 unidades", "Región de Ñuble", 83, 6),
    @(44616, "Primera", 200, 600, 700, 650, "$/paquete 6 unidades", "Región de Ñuble", 108, 6),
    @(44616, "Segunda", 100, 500, 500, 500, "$/paquete 6 unidades", "Región de Ñuble", 83, 6),
    @(44658, "Primera", 110, 6000, 7000, 6545, "$/paquete 36 unidades", "Región Metropolitana", 182, 36),
    @(44491, "Primera", 200, 600, 700, 650, "$/paquete 6 unidades", "Región Metropolitana", 108, 6),
    @(44491, "Segunda", 100, 500, 500, 500, "$/paquete 6 unidades", "Región Metropolitana", 83, 6),
    @(44631, "Primera", 220, 6000, 6500, 6227, "$/paquete 36 unidades", "Región Metropolitana", 173, 36),
    @(44637, "Primera", 110, 6500, 7000, 6773, "$/paquete 36 unidades", "Región Metropolitana", 188, 36),
    @(44649, "Primera", 220, 8000, 8500, 8227, "$/paquete 36 unidades", "Región Metropolitana", 229, 36),
    @(44554, "Primera", 200, 600, 700, 650, "$/paquete 6 unidades", "Región de Ñuble", 108, 6),
    @(44554, "Segunda", 100, 500, 500, 500, "$/paquete 6 unidades", "Región de Ñuble", 83, 6),
    @(44188, "Primera", 200, 600, 700, 650, "$/paquete 6 unidades", "Región de Ñuble", 108, 6),
    @(44188, "Segunda", 100, 500, 500, 500, "$/paquete 6 unidades", "Región de Ñuble", 83, 6),
    @(44358, "Primera", 200, 600, 700, 650, "$/paquete 6 unidades", "Región de Ñuble", 108, 6),
    @(44358, "Segunda", 100, 500, 500, 500, "$/paquete 6 unidades", "Región de Ñuble", 83, 6),
    @(44525, "Primera", 200, 600, 700, 650, "$/paquete 6 unidades", "Región de Ñuble", 108, 6),
    @(44525, "Segunda", 100, 500, 500, 500, "$/paquete 6 unidades", "Región de Ñuble", 83, 6),
    @(44644, "Primera", 160, 6500, 7000, 6750, "$/paquete 36 unidades", "Región Metropolitana", 188, 36),
    @(44335, "Primera", 150, 600, 700, 633, "$/paquete 6 unidades", "Región de Ñuble", 106, 6),
    @(44335, "Segunda", 50, 500, 500, 500, "$/paquete 6 unidades", "Región de Ñuble", 83, 6),
    @(44328, "Primera", 100, 600, 700, 650, "$/paquete 6 unidades", "Región de Ñuble", 108, 6),
    @(44328, "Segunda", 50, 500, 500, 500, "$/paquete 6 unidades", "Región de Ñuble", 83, 6),
    @(44643, "Primera", 180, 6500, 7000, 6778, "$/paquete 36 unidades", "Región Metropolitana", 188, 36),
    @(44659, "Primera", 300, 8000, 8500, 8250, "$/paquete 36 unidades", "Región Metropolitana", 229, 36)
)

$startRow = 2
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = $startRow + $i
    $data = $rowsData[$i]
    $ws.Cells.Item($r, 4).Value  = $data[0]   # D  Fecha
    $ws.Cells.Item($r, 9).Value  = $data[1]   # I  Calidad
    $ws.Cells.Item($r, 10).Value = $data[2]   # J  Volumen
    $ws.Cells.Item($r, 11).Value = $data[3]   # K  Precio minimo
    $ws.Cells.Item($r, 12).Value = $data[4]   # L  Precio maximo
    $ws.Cells.Item($r, 13).Value = $data[5]   # M  Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $data[6]   # N  Unidad de comercializacion
    $ws.Cells.Item($r, 15).Value = $data[7]   # O  Origen
    $ws.Cells.Item($r, 16).Value = $data[8]   # P  Precio $/Kg
    $ws.Cells.Item($r, 17).Value = $data[9]   # Q  Kg o Unidades
}
